$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("86-35=", $true, $false, $false, $false, $false, $true, 1, $false, "57-27=", 2)
$null = $d.Content.Find.Execute("45-18=", $true, $false, $false, $false, $false, $true, 1, $false, "4+11=", 2)
$null = $d.Content.Find.Execute("48-21=", $true, $false, $false, $false, $false, $true, 1, $false, "58+6=", 2)
$null = $d.Content.Find.Execute("21+42=", $true, $false, $false, $false, $false, $true, 1, $false, "40-40=", 2)
$null = $d.Content.Find.Execute("84-67=", $true, $false, $false, $false, $false, $true, 1, $false, "48+15=", 2)
$null = $d.Content.Find.Execute("19+78=", $true, $false, $false, $false, $false, $true, 1, $false, "86-38=", 2)
$null = $d.Content.Find.Execute("47-45=", $true, $false, $false, $false, $false, $true, 1, $false, "79-70=", 2)
$null = $d.Content.Find.Execute("59-50=", $true, $false, $false, $false, $false, $true, 1, $false, "11+29=", 2)
$null = $d.Content.Find.Execute("33+23=", $true, $false, $false, $false, $false, $true, 1, $false, "28+54=", 2)
$null = $d.Content.Find.Execute("35+55=", $true, $false, $false, $false, $false, $true, 1, $false, "61-6=", 2)
$null = $d.Content.Find.Execute("46+17=", $true, $false, $false, $false, $false, $true, 1, $false, "14+60=", 2)
$null = $d.Content.Find.Execute("75-44=", $true, $false, $false, $false, $false, $true, 1, $false, "80-58=", 2)
$null = $d.Content.Find.Execute("77-67=", $true, $false, $false, $false, $false, $true, 1, $false, "0+92=", 2)
$null = $d.Content.Find.Execute("87-76=", $true, $false, $false, $false, $false, $true, 1, $false, "51-34=", 2)
$null = $d.Content.Find.Execute("69-14=", $true, $false, $false, $false, $false, $true, 1, $false, "29+44=", 2)
$null = $d.Content.Find.Execute("51-15=", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=", 2)
$null = $d.Content.Find.Execute("33+52=", $true, $false, $false, $false, $false, $true, 1, $false, "92-35=", 2)
$null = $d.Content.Find.Execute("74-49=", $true, $false, $false, $false, $false, $true, 1, $false, "22+21=", 2)
$null = $d.Content.Find.Execute("0+80=", $true, $false, $false, $false, $false, $true, 1, $false, "12-8=", 2)
$null = $d.Content.Find.Execute("74-72=", $true, $false, $false, $false, $false, $true, 1, $false, "77-3=", 2)
$null = $d.Content.Find.Execute("66+14=", $true, $false, $false, $false, $false, $true, 1, $false, "95-63=", 2)
$null = $d.Content.Find.Execute("73-37=", $true, $false, $false, $false, $false, $true, 1, $false, "32+21=", 2)
$null = $d.Content.Find.Execute("19+45=", $true, $false, $false, $false, $false, $true, 1, $false, "80-73=", 2)
$null = $d.Content.Find.Execute("54-1=", $true, $false, $false, $false, $false, $true, 1, $false, "9+48=", 2)
$null = $d.Content.Find.Execute("43+52=", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=", 2)
$null = $d.Content.Find.Execute("11+11=", $true, $false, $false, $false, $false, $true, 1, $false, "26+34=", 2)
$null = $d.Content.Find.Execute("50-13=", $true, $false, $false, $false, $false, $true, 1, $false, "8+65=", 2)
$null = $d.Content.Find.Execute("86-31=", $true, $false, $false, $false, $false, $true, 1, $false, "81-0=", 2)
$null = $d.Content.Find.Execute("42+28=", $true, $false, $false, $false, $false, $true, 1, $false, "25+35=", 2)
$null = $d.Content.Find.Execute("61+10=", $true, $false, $false, $false, $false, $true, 1, $false, "12+58=", 2)
$null = $d.Content.Find.Execute("42+17=", $true, $false, $false, $false, $false, $true, 1, $false, "16+22=", 2)
$null = $d.Content.Find.Execute("82+9=", $true, $false, $false, $false, $false, $true, 1, $false, "8+78=", 2)
$null = $d.Content.Find.Execute("79-45=", $true, $false, $false, $false, $false, $true, 1, $false, "75-25=", 2)
$null = $d.Content.Find.Execute("51-13=", $true, $false, $false, $false, $false, $true, 1, $false, "36+41=", 2)
$null = $d.Content.Find.Execute("31+16=", $true, $false, $false, $false, $false, $true, 1, $false, "23+37=", 2)
$null = $d.Content.Find.Execute("77-19=", $true, $false, $false, $false, $false, $true, 1, $false, "77-24=", 2)
$null = $d.Content.Find.Execute("72-70=", $true, $false, $false, $false, $false, $true, 1, $false, "94+0=", 2)
$null = $d.Content.Find.Execute("72+3=", $true, $false, $false, $false, $false, $true, 1, $false, "87+1=", 2)
$null = $d.Content.Find.Execute("65+12=", $true, $false, $false, $false, $false, $true, 1, $false, "59-29=", 2)
$null = $d.Content.Find.Execute("55-4=", $true, $false, $false, $false, $false, $true, 1, $false, "93-72=", 2)
$null = $d.Content.Find.Execute("80-8=", $true, $false, $false, $false, $false, $true, 1, $false, "46-43=", 2)
$null = $d.Content.Find.Execute("56-30=", $true, $false, $false, $false, $false, $true, 1, $false, "45+48=", 2)
$null = $d.Content.Find.Execute("25+51=", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=", 2)
$null = $d.Content.Find.Execute("35-19=", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=", 2)
$null = $d.Content.Find.Execute("78-60=", $true, $false, $false, $false, $false, $true, 1, $false, "70-3=", 2)
$null = $d.Content.Find.Execute("84-80=", $true, $false, $false, $false, $false, $true, 1, $false, "14+71=", 2)
$null = $d.Content.Find.Execute("22+66=", $true, $false, $false, $false, $false, $true, 1, $false, "40-25=", 2)
$null = $d.Content.Find.Execute("53-34=", $true, $false, $false, $false, $false, $true, 1, $false, "49-47=", 2)
$null = $d.Content.Find.Execute("54-42=", $true, $false, $false, $false, $false, $true, 1, $false, "15+53=", 2)
$null = $d.Content.Find.Execute("59-42=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 2)
$null = $d.Content.Find.Execute("45+3=", $true, $false, $false, $false, $false, $true, 1, $false, "8-6=", 2)
$null = $d.Content.Find.Execute("27+55=", $true, $false, $false, $false, $false, $true, 1, $false, "81-34=", 2)
$null = $d.Content.Find.Execute("63-12=", $true, $false, $false, $false, $false, $true, 1, $false, "95+1=", 2)
$null = $d.Content.Find.Execute("90-78=", $true, $false, $false, $false, $false, $true, 1, $false, "19+75=", 2)
$null = $d.Content.Find.Execute("53-26=", $true, $false, $false, $false, $false, $true, 1, $false, "20+17=", 2)
$null = $d.Content.Find.Execute("68+26=", $true, $false, $false, $false, $false, $true, 1, $false, "25+68=", 2)
$null = $d.Content.Find.Execute("71-24=", $true, $false, $false, $false, $false, $true, 1, $false, "76+3=", 2)
$null = $d.Content.Find.Execute("31+65=", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=", 2)
$null = $d.Content.Find.Execute("76+5=", $true, $false, $false, $false, $false, $true, 1, $false, "40-6=", 2)
$null = $d.Content.Find.Execute("64-49=", $true, $false, $false, $false, $false, $true, 1, $false, "65+30=", 2)
$null = $d.Content.Find.Execute("81-8=", $true, $false, $false, $false, $false, $true, 1, $false, "39+51=", 2)
$null = $d.Content.Find.Execute("35+63=", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=", 2)
$null = $d.Content.Find.Execute("1-1=", $true, $false, $false, $false, $false, $true, 1, $false, "47-41=", 2)
$null = $d.Content.Find.Execute("0+76=", $true, $false, $false, $false, $false, $true, 1, $false, "9-0=", 2)
$null = $d.Content.Find.Execute("94-57=", $true, $false, $false, $false, $false, $true, 1, $false, "39+41=", 2)
$null = $d.Content.Find.Execute("91-83=", $true, $false, $false, $false, $false, $true, 1, $false, "71-13=", 2)
$null = $d.Content.Find.Execute("53-20=", $true, $false, $false, $false, $false, $true, 1, $false, "11+28=", 2)
$null = $d.Content.Find.Execute("18+9=", $true, $false, $false, $false, $false, $true, 1, $false, "38+40=", 2)
$null = $d.Content.Find.Execute("70-50=", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=", 2)
$null = $d.Content.Find.Execute("59-9=", $true, $false, $false, $false, $false, $true, 1, $false, "50-24=", 2)
$null = $d.Content.Find.Execute("23-14=", $true, $false, $false, $false, $false, $true, 1, $false, "87-50=", 2)
$null = $d.Content.Find.Execute("89-72=", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=", 2)
$null = $d.Content.Find.Execute("74-4=", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=", 2)
$null = $d.Content.Find.Execute("12+27=", $true, $false, $false, $false, $false, $true, 1, $false, "83+5=", 2)
$null = $d.Content.Find.Execute("25+70=", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=", 2)
$null = $d.Content.Find.Execute("79-1=", $true, $false, $false, $false, $false, $true, 1, $false, "58-48=", 2)
$null = $d.Content.Find.Execute("88-82=", $true, $false, $false, $false, $false, $true, 1, $false, "22+56=", 2)
$null = $d.Content.Find.Execute("4+34=", $true, $false, $false, $false, $false, $true, 1, $false, "28+48=", 2)
$null = $d.Content.Find.Execute("5+78=", $true, $false, $false, $false, $false, $true, 1, $false, "37-1=", 2)
$null = $d.Content.Find.Execute("99-58=", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=", 2)
$null = $d.Content.Find.Execute("38+20=", $true, $false, $false, $false, $false, $true, 1, $false, "62-8=", 2)
$null = $d.Content.Find.Execute("2+29=", $true, $false, $false, $false, $false, $true, 1, $false, "2+73=", 2)
$null = $d.Content.Find.Execute("52-52=", $true, $false, $false, $false, $false, $true, 1, $false, "70-54=", 2)
$null = $d.Content.Find.Execute("28+57=", $true, $false, $false, $false, $false, $true, 1, $false, "91-12=", 2)
$null = $d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=", 2)
$null = $d.Content.Find.Execute("47+43=", $true, $false, $false, $false, $false, $true, 1, $false, "96-54=", 2)
$null = $d.Content.Find.Execute("13+67=", $true, $false, $false, $false, $false, $true, 1, $false, "87-36=", 2)
$null = $d.Content.Find.Execute("58-11=", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=", 2)
$null = $d.Content.Find.Execute("60-5=", $true, $false, $false, $false, $false, $true, 1, $false, "36+59=", 2)
$null = $d.Content.Find.Execute("38-20=", $true, $false, $false, $false, $false, $true, 1, $false, "18+65=", 2)
$null = $d.Content.Find.Execute("64-1=", $true, $false, $false, $false, $false, $true, 1, $false, "76-64=", 2)
$null = $d.Content.Find.Execute("6+8=", $true, $false, $false, $false, $false, $true, 1, $false, "50-32=", 2)
$null = $d.Content.Find.Execute("34+31=", $true, $false, $false, $false, $false, $true, 1, $false, "9+9=", 2)
$null = $d.Content.Find.Execute("91-26=", $true, $false, $false, $false, $false, $true, 1, $false, "84-13=", 2)
$null = $d.Content.Find.Execute("33+2=", $true, $false, $false, $false, $false, $true, 1, $false, "0+91=", 2)
$null = $d.Content.Find.Execute("94-41=", $true, $false, $false, $false, $false, $true, 1, $false, "81-1=", 2)
$null = $d.Content.Find.Execute("40+10=", $true, $false, $false, $false, $false, $true, 1, $false, "86-46=", 2)
$null = $d.Content.Find.Execute("15+18=", $true, $false, $false, $false, $false, $true, 1, $false, "6+15=", 2)
$null = $d.Content.Find.Execute("66-63=", $true, $false, $false, $false, $false, $true, 1, $false, "97-49=", 2)
$null = $d.Content.Find.Execute("18+16=", $true, $false, $false, $false, $false, $true, 1, $false, "72-11=", 2)
